$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The real edit: F3's hard-coded value changed (954859 -> 945169).
# F4 (=F2-F3) and F5 (=F4/F2) are formulas and recalc automatically.
$ws.Range("F3").Value = 945169

# Match the final selection state as closely as the object model allows.
# The authored change left a multi-area selection (F3 and F13, active cell
# F13); select both areas and make F13 the active cell.
$ws.Range("F3,F13").Select()
$ws.Range("F13").Activate()
